$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.506.13'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').Value = '2.223.96'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '270.10'
$ws.Range('E5').Value = '  +3.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '92.19'
$ws.Range('E6').Value = '  +11.64%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').Value = '  -1.17%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '45.67'
$ws.Range('E10').Value = '  +4.12%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0927'
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.22'
$ws.Range('E12').Value = '  +16.46%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.104'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('D14').Value = '2.562.02'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.03'
$ws.Range('E15').Value = '  +3.23%  '
$ws.Range('D16').Value = '2.229.60'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.800'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('D18').Value = '43.489.06'
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.98'
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '70.27'
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '232.46'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.02'
$ws.Range('E24').Value = '  -3.41%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.50'
$ws.Range('E26').Value = '  +11.12%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.27'
$ws.Range('E27').Value = '  +4.38%  '
$ws.Range('E28').Value = '  +5.33%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '40.53'
$ws.Range('E29').Value = '  -2.27%  '
$ws.Range('E30').Value = '  +1.85%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '172.42'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0919'
$ws.Range('E32').Value = '  +4.91%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.79'
$ws.Range('E33').Value = '  +0.67%  '
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.111'
$ws.Range('E36').Value = '  -3.97%  '
$ws.Range('E37').Value = '  -3.02%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.29'
$ws.Range('E38').Value = '  -4.80%  '
$ws.Range('E39').Value = '  +20.28%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '12.44'
$ws.Range('E40').Value = '  -8.50%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.16'
$ws.Range('E41').Value = '  +1.97%  '
$ws.Range('E42').Value = '  +8.06%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '63.14'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.31'
$ws.Range('E44').Value = '  -4.74%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0986'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '8.35'
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '100.10'
$ws.Range('E47').Value = '  -2.91%  '
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.437'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.46'
$ws.Range('E51').Value = '  -7.24%  '
